# Apply updated Fecha / Volumen / Precio mínimo / Precio máximo /
# Precio promedio ponderado / Precio $/Kg values (weekly refresh of the
# "Hortaliza, Agrícola del Norte S.A. de Arica - Jengibre" sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row -> @(Fecha(D), Volumen(J), PrecioMinimo(K), PrecioMaximo(L), PrecioPromedioPonderado(M), PrecioKg(P))
$rows = @{
    2  = @(44832, 100, 13000, 14000, 13500, 1038)
    3  = @(44616, 120, 19000, 20000, 19500, 1500)
    4  = @(44469, 140, 13000, 14000, 13500, 1038)
    5  = @(44379, 120, 12000, 13000, 12667,  974)
    6  = @(44580, 160, 11000, 12000, 11500,  885)
    8  = @(44320, 160, 19000, 20000, 19500, 1500)
    9  = @(44855, 500, 10000, 10000, 10000,  769)
    10 = @(44389, 120, 12000, 13000, 12500,  962)
    11 = @(44764, 200, 12000, 13000, 12500,  962)
    12 = @(44397, 140, 12500, 13000, 12750,  981)
    13 = @(44592, 120, 12000, 13000, 12500,  962)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals[0]
    $ws.Range("J$r").Value = $vals[1]
    $ws.Range("K$r").Value = $vals[2]
    $ws.Range("L$r").Value = $vals[3]
    $ws.Range("M$r").Value = $vals[4]
    $ws.Range("P$r").Value = $vals[5]
}
